$wb = $excel.ActiveWorkbook

$wsPalavras = $wb.Worksheets.Item("PalavrasReservadas")
$wsSinais = $wb.Worksheets.Item("Sinais")
$wsTransicoes = $wb.Worksheets.Item("Transicoes")

# Add new reserved words / signs in the exact order they were first typed so
# that the shared-strings table receives the same index assignment as the
# original edit.

# 1) function / main / return on PalavrasReservadas
$wsPalavras.Range("A10").Value = 9
$wsPalavras.Range("B10").Value = "function"

$wsPalavras.Range("A11").Value = 10
$wsPalavras.Range("B11").Value = "main"

$wsPalavras.Range("A12").Value = 11
$wsPalavras.Range("B12").Value = "return"

# 2) { and } on Sinais (rows 14-15)
$wsSinais.Range("A14").Value = 13
$wsSinais.Range("B14").Value = "{"

$wsSinais.Range("A15").Value = 14
$wsSinais.Range("B15").Value = "}"

# 3) and / or / not on PalavrasReservadas (rows 13-15)
$wsPalavras.Range("A13").Value = 12
$wsPalavras.Range("B13").Value = "and"

$wsPalavras.Range("A14").Value = 13
$wsPalavras.Range("B14").Value = "or"

$wsPalavras.Range("A15").Value = 14
$wsPalavras.Range("B15").Value = "not"

# 4) ! on Sinais (row 13)
$wsSinais.Range("A13").Value = 12
$wsSinais.Range("B13").Value = "!"

# Update selections on each sheet to match the final workbook state.
$wsPalavras.Activate()
$wsPalavras.Range("A16").Select()

$wsTransicoes.Activate()
$wsTransicoes.Range("C11").Select()

# Sinais ends up the active (visible) tab with B16 selected.
$wsSinais.Activate()
$wsSinais.Range("B16").Select()
